# Auto commit at 2025-12-09  8:05:40.75
# Appends two new data rows (198 and 199) to Sheet1, mirroring the
# existing pattern of paired "四方坪站充电量(kw)" / "高岭站充电量(kw)" rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 198: 四方坪站充电量(kw) for 2025-12-08 (serial 45999)
$ws.Range("A198").Value = 45999
$ws.Range("A198").NumberFormat = "yyyy\-mm\-dd"
$ws.Range("B198").Value = "四方坪站充电量(kw)"

$row198 = @(612.94600000000003,1315.7599999999998,342.64100000000002,145.89500000000001,267.75899999999996,539.803,458.00300000000004,86.198000000000008,66.099999999999994,193.352,222.96999999999997,365.62799999999999,926.07399999999996,1568.4410000000003,439.97999999999996,393.37500000000006,251.15800000000004,220.875,43.620000000000005,128.489,0,93.503000000000014,67.8,89.000999999999991)

for ($i = 0; $i -lt $row198.Length; $i++) {
    $col = 3 + $i
    $cell = $ws.Cells.Item(198, $col)
    $cell.Value = $row198[$i]
    $cell.NumberFormat = "0.00_);[Red]\(0.00\)"
}

# Row 199: 高岭站充电量(kw) for 2025-12-08 (serial 45999)
$ws.Range("A199").Value = 45999
$ws.Range("A199").NumberFormat = "yyyy\-mm\-dd"
$ws.Range("B199").Value = "高岭站充电量(kw)"

$row199 = @(412.91799999999995,831.49,0,110.71900000000001,43.258000000000003,117.946,133.411,216.82700000000003,291.89499999999998,194.64199999999997,143.73699999999999,230.38,538.95100000000002,378.39400000000006,426.08199999999994,201.82999999999998,53.960999999999999,132.50799999999998,59.320999999999998,130.00200000000001,96.117000000000004,121.33199999999999,0,50.734999999999999)

for ($i = 0; $i -lt $row199.Length; $i++) {
    $col = 3 + $i
    $cell = $ws.Cells.Item(199, $col)
    $cell.Value = $row199[$i]
    $cell.NumberFormat = "0.00_);[Red]\(0.00\)"
}

# Update the saved view state to match the committed worksheet view.
$ws.Range("H208").Select()
